# Apply the BiblicaStudyNotesKeyTerms resource-data update.
#
# Before, the paragraphs (1-based) relevant to this edit look like:
#   3  Heading2  "Resource: 關鍵詞 (Biblica)"
#   4  Heading2  "License Information"                              <- remove entirely
#   5  Normal    "關鍵詞 (Biblica) (Chinese (Traditional)) is based  <- rewrite
#               on: Biblica Bible Dictionary, Biblica, Inc., 2023,
#               which is licensed under a CC BY-SA 4.0 license."
#   6  Normal    "This PDF version is provided under the same        <- remove entirely
#               license."
#  ...
#  11  Normal    "憐憫" (italic)                                     <- remove entirely
#
# Work bottom-to-top for the paragraph deletions so earlier indices stay valid,
# and do the in-place text edits of paragraph 5 first (that one keeps its index).

$d = $word.ActiveDocument

# Guard rails: make sure the paragraph indices below still point at the
# paragraphs we expect before mutating anything.
if ($d.Paragraphs.Item(4).Range.Text.Trim() -ne "License Information") {
    throw "paragraph 4 is not 'License Information' - aborting"
}
if ($d.Paragraphs.Item(6).Range.Text.Trim() -ne "This PDF version is provided under the same license.") {
    throw "paragraph 6 is not the PDF-license sentence - aborting"
}
if ($d.Paragraphs.Item(11).Range.Text.Trim() -ne "憐憫") {
    throw "paragraph 11 is not the stand-alone '憐憫' paragraph - aborting"
}

# ---------------------------------------------------------------------------
# 1) Paragraph 5: rewrite the "關鍵詞 (Biblica) ... license." blurb into the
#    new "Biblica Study Notes (Key Terms) ... Mission Mutual." blurb.
# ---------------------------------------------------------------------------

# 1a) Drop the dictionary credit + both hyperlinked citations + trailing
#     period; this also removes the two hyperlinks from the paragraph.
$rng = $d.Paragraphs.Item(5).Range
$rng.Find.Execute(
    "Biblica Bible Dictionary, Biblica, Inc., 2023, which is licensed under a CC BY-SA 4.0 license.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "", 2) | Out-Null

# 1b) Bold run: "關鍵詞 (Biblica)" -> "Biblica Study Notes (Key Terms)"
$rng = $d.Paragraphs.Item(5).Range
$rng.Find.Execute(
    "關鍵詞 (Biblica)",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Biblica Study Notes (Key Terms)", 2) | Out-Null

# 1c) Replace the remaining tail (" (Chinese (Traditional)) is based on: ")
#     with the new license blurb + adaptation-languages sentence.
$rng = $d.Paragraphs.Item(5).Range
$rng.Find.Execute(
    " (Chinese (Traditional)) is based on: ",
    $false, $false, $false, $false, $false, $true, 1, $false,
    " © 2023 Biblica Inc. Released under CC BY-SA 4.0 license. Biblica Study Notes has been adapted in the following languages: Tok Pisin, Arabic (عربي), French (Français), Hindi (हिंदी), Indonesian (Bahasa Indonesia), Portuguese (Português), Russian (Русский), Spanish (Español), Swahili (Kiswahili), and Simplified Chinese (简体中文)from Biblica Study Notes © 2023 Biblica Inc. Released under CC BY-SA 4.0 license by Mission Mutual.",
    2) | Out-Null

# ---------------------------------------------------------------------------
# 2) Remove paragraph 6, "This PDF version is provided under the same
#    license." (index unaffected by the paragraph-5 text edits above).
# ---------------------------------------------------------------------------
$d.Paragraphs.Item(6).Range.Delete()

# ---------------------------------------------------------------------------
# 3) Remove paragraph 4, the "License Information" Heading2.
# ---------------------------------------------------------------------------
$d.Paragraphs.Item(4).Range.Delete()

# ---------------------------------------------------------------------------
# 4) Remove the stand-alone italic "憐憫" paragraph. It was paragraph 11
#    originally; after deleting the two paragraphs above it, it is now at 9.
# ---------------------------------------------------------------------------
if ($d.Paragraphs.Item(9).Range.Text.Trim() -ne "憐憫") {
    throw "paragraph 9 (post-delete) is not the stand-alone '憐憫' paragraph - aborting"
}
$d.Paragraphs.Item(9).Range.Delete()
